$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
# E1 text changes from "(f-f')^2" to "(f-f')^2 normalized"
$ws.Range("E1").Value = "(f-f')^2 normalized"
# New F1 header "(f-f')^2 classic", formatted like the other header cells
$ws.Range("F1").Value = "(f-f')^2 classic"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# --- Extend A column (rows 10-22 => x0 = 8..20) with same formatting as existing A cells ---
$ws.Range("A2").Copy()
$ws.Range("A10:A22").PasteSpecial(-4122)

# --- Data rows ---
# A=x0, B=f, C=f-prime, D=f-prime (regressed), E=(f-f-prime)^2 normalized, F=(f-f-prime)^2 classic
$ws.Range("B2").Value = 83521
$ws.Range("C2").Value = 87158
$ws.Range("D2").Value = [double]"87157.99597575798"
$ws.Range("E2").Value = [double]"2.655876616128772E-14"
$ws.Range("F2").Value = [double]"1.619452385035841E-05"
$ws.Range("B3").Value = 81632
$ws.Range("C3").Value = 74287
$ws.Range("D3").Value = [double]"74286.99524180366"
$ws.Range("E3").Value = [double]"3.712995541860694E-14"
$ws.Range("F3").Value = [double]"2.264043245264865E-05"
$ws.Range("B4").Value = 82277
$ws.Range("C4").Value = 62375
$ws.Range("D4").Value = [double]"62374.96224954988"
$ws.Range("E4").Value = [double]"2.337135969570663E-12"
$ws.Range("F4").Value = [double]"0.001425096483938641"
$ws.Range("B5").Value = 78349
$ws.Range("C5").Value = 63904
$ws.Range("D5").Value = [double]"63904.00627640509"
$ws.Range("E5").Value = [double]"6.46043323492143E-14"
$ws.Range("F5").Value = [double]"3.939326084136255E-05"
$ws.Range("B6").Value = 78216
$ws.Range("C6").Value = 61954
$ws.Range("D6").Value = [double]"61954.00531282049"
$ws.Range("E6").Value = [double]"4.629029997574734E-14"
$ws.Range("F6").Value = [double]"2.822606154889799E-05"
$ws.Range("B7").Value = 63483
$ws.Range("C7").Value = 62836
$ws.Range("D7").Value = [double]"62835.9974089648"
$ws.Range("E7").Value = [double]"1.1009975084499E-14"
$ws.Range("F7").Value = [double]"6.71346339643649E-06"
$ws.Range("B8").Value = 52243
$ws.Range("C8").Value = 67016
$ws.Range("D8").Value = [double]"67016.009753197"
$ws.Range("E8").Value = [double]"1.560032707885569E-13"
$ws.Range("F8").Value = [double]"9.512485162341112E-05"
$ws.Range("B9").Value = 53558
$ws.Range("C9").Value = 64482
$ws.Range("D9").Value = [double]"64482.00568542173"
$ws.Range("E9").Value = [double]"5.301088809124646E-14"
$ws.Range("F9").Value = [double]"3.232402019681284E-05"
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 54685
$ws.Range("C10").Value = 62172
$ws.Range("D10").Value = [double]"62172.00532368233"
$ws.Range("E10").Value = [double]"4.647977059794016E-14"
$ws.Range("F10").Value = [double]"2.834159354668098E-05"
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 54507
$ws.Range("C11").Value = 49524
$ws.Range("D11").Value = [double]"49523.99757949099"
$ws.Range("E11").Value = [double]"9.608445204132538E-15"
$ws.Range("F11").Value = [double]"5.858863860307554E-06"
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 48757
$ws.Range("C12").Value = 48480
$ws.Range("D12").Value = [double]"48480.00462402251"
$ws.Range("E12").Value = [double]"3.506546409980269E-14"
$ws.Range("F12").Value = [double]"2.138158419954067E-05"
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 44455
$ws.Range("C13").Value = 44876
$ws.Range("D13").Value = [double]"44876.00417088519"
$ws.Range("E13").Value = [double]"2.852963272457853E-14"
$ws.Range("F13").Value = [double]"1.739628323384087E-05"
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 48869
$ws.Range("C14").Value = 31926
$ws.Range("D14").Value = [double]"31926.00179487233"
$ws.Range("E14").Value = [double]"5.283319024730905E-15"
$ws.Range("F14").Value = [double]"3.221566677054374E-06"
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 32235
$ws.Range("C15").Value = 42462
$ws.Range("D15").Value = [double]"42461.99784561976"
$ws.Range("E15").Value = [double]"7.61174844261693E-15"
$ws.Range("F15").Value = [double]"4.641354232097346E-06"
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 41006
$ws.Range("C16").Value = 21787
$ws.Range("D16").Value = [double]"21787.00021050921"
$ws.Range("E16").Value = [double]"7.267447477178973E-17"
$ws.Range("F16").Value = [double]"4.431412590613688E-08"
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 35150
$ws.Range("C17").Value = 20550
$ws.Range("D17").Value = [double]"20550.01211102101"
$ws.Range("E17").Value = [double]"2.405477098076482E-13"
$ws.Range("F17").Value = [double]"0.0001466768300647205"
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 29272
$ws.Range("C18").Value = 15917
$ws.Range("D18").Value = [double]"15916.99429901746"
$ws.Range("E18").Value = [double]"5.330146328441325E-14"
$ws.Range("F18").Value = [double]"3.250120187494834E-05"
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 28595
$ws.Range("C19").Value = 15775
$ws.Range("D19").Value = [double]"15775.0037127686"
$ws.Range("E19").Value = [double]"2.260661180455557E-14"
$ws.Range("F19").Value = [double]"1.378465071130766E-05"
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 15816
$ws.Range("C20").Value = 14717
$ws.Range("D20").Value = [double]"14717.00456064919"
$ws.Range("E20").Value = [double]"3.411088959443608E-14"
$ws.Range("F20").Value = [double]"2.079952104576056E-05"
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 13053
$ws.Range("C21").Value = 13947
$ws.Range("D21").Value = [double]"13947.00721488086"
$ws.Range("E21").Value = [double]"8.536857640679696E-14"
$ws.Range("F21").Value = [double]"5.205450584251848E-05"
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 14595
$ws.Range("C22").Value = 11970
$ws.Range("D22").Value = [double]"11969.99879115759"
$ws.Range("E22").Value = [double]"2.396509134882938E-15"
$ws.Range("F22").Value = [double]"1.46129997673601E-06"
